$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H95").Value = 25000
$ws.Range("J95").Value = 25000
$ws.Range("L95").Value = 25000
$ws.Range("N95").Value = -30492

$ws.Range("H137").Value = 44219.645
$ws.Range("I137").Value = 55343.727
$ws.Range("J137").Value = 3431.3333
$ws.Range("K137").Value = 166031.181
$ws.Range("L137").Value = 10293.9999
$ws.Range("M137").Value = -163481.181
$ws.Range("N137").Value = -15393.9999

$ws.Range("H141").Value = 10983.857
$ws.Range("I141").Value = 10983.857
$ws.Range("K141").Value = 32951.571
$ws.Range("M141").Value = -27771.571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 8912.666999999999
$ws.Range("I61").Value = 10871.75
$ws.Range("K61").Value = 10871.75
$ws.Range("M61").Value = -10659.75

$ws.Range("H97").Value = 792095.25
$ws.Range("I97").Value = 1011491.6
$ws.Range("K97").Value = 1011491.6
$ws.Range("M97").Value = -1010995.6

$ws.Range("H136").Value = 8912.666999999999
$ws.Range("I136").Value = 10871.75
$ws.Range("K136").Value = 32615.25
$ws.Range("M136").Value = -30065.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2781258.5
$ws.Range("I94").Value = 4348472
$ws.Range("J94").Value = 8496.154
$ws.Range("K94").Value = 4348472
$ws.Range("L94").Value = 8496.154
$ws.Range("M94").Value = -4348021
$ws.Range("N94").Value = -9398.154

$ws.Range("H134").Value = 3547.122
$ws.Range("I134").Value = 1718.2916
$ws.Range("J134").Value = 6129
$ws.Range("K134").Value = 5154.8748
$ws.Range("L134").Value = 18387
$ws.Range("M134").Value = -2619.8748
$ws.Range("N134").Value = -23457

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 20010.5
$ws.Range("J29").Value = 20010.5
$ws.Range("L29").Value = 20010.5
$ws.Range("N29").Value = -20596.5

$ws.Range("H58").Value = 6172.846
$ws.Range("I58").Value = 7793.5625
$ws.Range("K58").Value = 7793.5625
$ws.Range("M58").Value = -7590.5625

$ws.Range("H104").Value = 64578
$ws.Range("J104").Value = 64578
$ws.Range("L104").Value = 64578
$ws.Range("N104").Value = -69820

$ws.Range("H105").Value = 1968.1428
$ws.Range("I105").Value = 1299
$ws.Range("K105").Value = 1299
$ws.Range("M105").Value = 448

$ws.Range("H117").Value = 64995
$ws.Range("J117").Value = 64995
$ws.Range("L117").Value = 64995
$ws.Range("N117").Value = -74173

$ws.Range("H136").Value = 6172.846
$ws.Range("I136").Value = 7793.5625
$ws.Range("K136").Value = 23380.6875
$ws.Range("M136").Value = -20830.6875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 259.625
$ws.Range("J2").Value = 343.7619
$ws.Range("L2").Value = 2062.5714
$ws.Range("N2").Value = -2288.5714

$ws.Range("H38").Value = 140.90909
$ws.Range("I38").Value = 71
$ws.Range("K38").Value = 213
$ws.Range("M38").Value = 134

$ws.Range("H98").Value = 1677.5714
$ws.Range("J98").Value = 1845.5454
$ws.Range("L98").Value = 5536.6362
$ws.Range("N98").Value = -8532.636200000001

$ws.Range("H114").Value = 836157
$ws.Range("J114").Value = 912168.75
$ws.Range("L114").Value = 2736506.25
$ws.Range("N114").Value = -2743014.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 10027
$ws.Range("I31").Value = 5040.5
$ws.Range("J31").Value = 20000
$ws.Range("K31").Value = 5040.5
$ws.Range("L31").Value = 20000
$ws.Range("M31").Value = -4748.5
$ws.Range("N31").Value = -20584

$ws.Range("H37").Value = 10027
$ws.Range("I37").Value = 5040.5
$ws.Range("J37").Value = 20000
$ws.Range("K37").Value = 5040.5
$ws.Range("L37").Value = 20000
$ws.Range("M37").Value = -4763.5
$ws.Range("N37").Value = -20554

$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H95").Value = 37597.332
$ws.Range("J95").Value = 37597.332
$ws.Range("L95").Value = 37597.332
$ws.Range("N95").Value = -43089.332

$ws.Range("H132").Value = 3414.5557
$ws.Range("I132").Value = 3100.4092
$ws.Range("J132").Value = 4796.8
$ws.Range("K132").Value = 9301.2276
$ws.Range("L132").Value = 14390.4
$ws.Range("M132").Value = -6771.2276
$ws.Range("N132").Value = -19450.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 10000
$ws.Range("I25").Value = 10000
$ws.Range("K25").Value = 10000
$ws.Range("M25").Value = -9770

$ws.Range("H61").Value = 37042852
$ws.Range("I61").Value = 55556270
$ws.Range("K61").Value = 55556270
$ws.Range("M61").Value = -55556068

$ws.Range("H68").Value = 1999
$ws.Range("I68").Value = 1999
$ws.Range("K68").Value = 1999
$ws.Range("M68").Value = -1250

$ws.Range("H71").Value = 1999
$ws.Range("I71").Value = 1999
$ws.Range("K71").Value = 9995
$ws.Range("M71").Value = -6251

$ws.Range("H93").Value = 18530280
$ws.Range("I93").Value = 23810224
$ws.Range("J93").Value = 50474.75
$ws.Range("K93").Value = 23810224
$ws.Range("L93").Value = 50474.75
$ws.Range("M93").Value = -23808976
$ws.Range("N93").Value = -52970.75

$ws.Range("H95").Value = 21333
$ws.Range("J95").Value = 21333
$ws.Range("L95").Value = 21333
$ws.Range("N95").Value = -26825

$ws.Range("H97").Value = 54047
$ws.Range("J97").Value = 54047
$ws.Range("L97").Value = 54047
$ws.Range("N97").Value = -56029

$ws.Range("H113").Value = 37042852
$ws.Range("I113").Value = 55556270
$ws.Range("K113").Value = 55556270
$ws.Range("M113").Value = -55554100

$ws.Range("H136").Value = 31526.553
$ws.Range("I136").Value = 47091.695
$ws.Range("K136").Value = 141275.085
$ws.Range("M136").Value = -138725.085

$ws.Range("H140").Value = 98089
$ws.Range("J140").Value = 98088.5
$ws.Range("L140").Value = 98088.5
$ws.Range("N140").Value = -108448.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 63849.4
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 63849.4
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 63849.4
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -64311.4

$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

$ws.Range("H134").Value = 63849.4
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 63849.4
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 191548.2
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -196618.2
